$d = $word.ActiveDocument

$replacements = @(
    @("711÷6=", "575÷5="),
    @("313÷5=", "442÷3="),
    @("395÷4=", "869÷6="),
    @("308÷9=", "491÷3="),
    @("327÷7=", "908÷2="),
    @("950÷9=", "679÷4="),
    @("627÷2=", "148÷9="),
    @("536÷7=", "180÷8="),
    @("748÷7=", "627÷7="),
    @("371÷8=", "227÷4="),
    @("568÷7=", "458÷2="),
    @("708÷2=", "894÷9="),
    @("439÷5=", "915÷5="),
    @("449÷2=", "887÷9="),
    @("407÷3=", "675÷5="),
    @("250÷4=", "782÷3="),
    @("842÷3=", "270÷5="),
    @("823÷8=", "759÷4="),
    @("543÷2=", "742÷5="),
    @("439÷2=", "714÷8="),
    @("768÷8=", "373÷2="),
    @("427÷8=", "885÷3="),
    @("280÷8=", "124÷6="),
    @("455÷6=", "605÷3="),
    @("944÷8=", "168÷7=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
